# Auto-generated script applying the betexplorer czech-republic fnl 2023-2024 update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shuffled match rows (home/away/odds/url); index/meta columns A-E unchanged ---

# Rows 19, 20, 21, 25 (rotate match data)
# Row 19
$ws.Range("F19").Value = "Varnsdorf"
$ws.Range("G19").Value = 3
$ws.Range("H19").Value = "Vlasim"
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 2.04
$ws.Range("K19").Value = "01/08/2023 06:12"
$ws.Range("L19").Value = 2.41
$ws.Range("M19").Value = "02/08/2023 17:56"
$ws.Range("N19").Value = 3.64
$ws.Range("O19").Value = "01/08/2023 06:12"
$ws.Range("P19").Value = 3.92
$ws.Range("Q19").Value = "02/08/2023 17:50"
$ws.Range("R19").Value = 3.07
$ws.Range("S19").Value = "01/08/2023 06:12"
$ws.Range("T19").Value = 2.58
$ws.Range("U19").Value = "02/08/2023 17:56"
$ws.Range("V19").Value = "https://www.betexplorer.com/football/czech-republic/fnl/varnsdorf-vlasim/pKQ0bRV3/"

# Row 20
$ws.Range("F20").Value = "Pribram"
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = "Opava"
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 2.14
$ws.Range("K20").Value = "01/08/2023 06:12"
$ws.Range("L20").Value = 2.39
$ws.Range("M20").Value = "02/08/2023 17:58"
$ws.Range("N20").Value = 3.41
$ws.Range("O20").Value = "01/08/2023 06:12"
$ws.Range("P20").Value = 3.42
$ws.Range("Q20").Value = "02/08/2023 17:58"
$ws.Range("R20").Value = 3.03
$ws.Range("S20").Value = "01/08/2023 06:12"
$ws.Range("T20").Value = 2.89
$ws.Range("U20").Value = "02/08/2023 17:58"
$ws.Range("V20").Value = "https://www.betexplorer.com/football/czech-republic/fnl/pribram-opava/xE3Zs6oN/"

# Row 21
$ws.Range("F21").Value = "Vyskov"
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = "Kromeriz"
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = 1.62
$ws.Range("K21").Value = "01/08/2023 13:11"
$ws.Range("L21").Value = 1.52
$ws.Range("M21").Value = "02/08/2023 17:52"
$ws.Range("N21").Value = 3.8
$ws.Range("O21").Value = "01/08/2023 13:11"
$ws.Range("P21").Value = 4.25
$ws.Range("Q21").Value = "02/08/2023 17:54"
$ws.Range("R21").Value = 4.52
$ws.Range("S21").Value = "01/08/2023 13:11"
$ws.Range("T21").Value = 6.12
$ws.Range("U21").Value = "02/08/2023 17:54"
$ws.Range("V21").Value = "https://www.betexplorer.com/football/czech-republic/fnl/mfk-vyskov-kromeriz/2eyda7Gc/"

# Row 25
$ws.Range("F25").Value = "Zizkov"
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = "Lisen"
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = 2.97
$ws.Range("K25").Value = "01/08/2023 13:11"
$ws.Range("L25").Value = 2.25
$ws.Range("M25").Value = "02/08/2023 17:55"
$ws.Range("N25").Value = 3.17
$ws.Range("O25").Value = "01/08/2023 13:11"
$ws.Range("P25").Value = 3.4
$ws.Range("Q25").Value = "02/08/2023 17:55"
$ws.Range("R25").Value = 2.23
$ws.Range("S25").Value = "01/08/2023 13:11"
$ws.Range("T25").Value = 3.13
$ws.Range("U25").Value = "02/08/2023 17:55"
$ws.Range("V25").Value = "https://www.betexplorer.com/football/czech-republic/fnl/zizkov-lisen/G42wsQ0T/"

# Rows 83, 84, 85 (rotate match data)
# Row 83
$ws.Range("F83").Value = "Opava"
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = "Brno"
$ws.Range("I83").Value = 1
$ws.Range("J83").Value = 2.86
$ws.Range("K83").Value = "28/09/2023 08:12"
$ws.Range("L83").Value = 2.06
$ws.Range("M83").Value = "30/09/2023 15:58"
$ws.Range("N83").Value = 3.16
$ws.Range("O83").Value = "28/09/2023 08:12"
$ws.Range("P83").Value = 3.45
$ws.Range("Q83").Value = "30/09/2023 15:59"
$ws.Range("R83").Value = 2.3
$ws.Range("S83").Value = "28/09/2023 08:12"
$ws.Range("T83").Value = 3.53
$ws.Range("U83").Value = "30/09/2023 15:58"
$ws.Range("V83").Value = "https://www.betexplorer.com/football/czech-republic/fnl/opava-brno/t6rMuyxb/"

# Row 84
$ws.Range("F84").Value = "Kromeriz"
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = "Vlasim"
$ws.Range("I84").Value = 1
$ws.Range("J84").Value = 2.67
$ws.Range("K84").Value = "28/09/2023 08:12"
$ws.Range("L84").Value = 3.1
$ws.Range("M84").Value = "30/09/2023 15:47"
$ws.Range("N84").Value = 3.28
$ws.Range("O84").Value = "28/09/2023 08:12"
$ws.Range("P84").Value = 3.57
$ws.Range("Q84").Value = "30/09/2023 15:47"
$ws.Range("R84").Value = 2.37
$ws.Range("S84").Value = "28/09/2023 08:12"
$ws.Range("T84").Value = 2.2
$ws.Range("U84").Value = "30/09/2023 15:47"
$ws.Range("V84").Value = "https://www.betexplorer.com/football/czech-republic/fnl/kromeriz-vlasim/lWoUwF6A/"

# Row 85
$ws.Range("F85").Value = "Lisen"
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = "Taborsko"
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1.95
$ws.Range("K85").Value = "28/09/2023 08:12"
$ws.Range("L85").Value = 2.16
$ws.Range("M85").Value = "30/09/2023 15:48"
$ws.Range("N85").Value = 3.3
$ws.Range("O85").Value = "28/09/2023 08:12"
$ws.Range("P85").Value = 3.31
$ws.Range("Q85").Value = "30/09/2023 15:58"
$ws.Range("R85").Value = 3.48
$ws.Range("S85").Value = "28/09/2023 08:12"
$ws.Range("T85").Value = 3.42
$ws.Range("U85").Value = "30/09/2023 15:48"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/czech-republic/fnl/lisen-taborsko/KdsQvei4/"

# Rows 94, 95 (swap match data)
# Row 94
$ws.Range("F94").Value = "Taborsko"
$ws.Range("G94").Value = 6
$ws.Range("H94").Value = "Sigma Olomouc B"
$ws.Range("I94").Value = 2
$ws.Range("J94").Value = 1.85
$ws.Range("K94").Value = "05/10/2023 08:12"
$ws.Range("L94").Value = 1.75
$ws.Range("M94").Value = "07/10/2023 15:03"
$ws.Range("N94").Value = 3.45
$ws.Range("O94").Value = "05/10/2023 08:12"
$ws.Range("P94").Value = 3.72
$ws.Range("Q94").Value = "07/10/2023 15:03"
$ws.Range("R94").Value = 3.65
$ws.Range("S94").Value = "05/10/2023 08:12"
$ws.Range("T94").Value = 4.63
$ws.Range("U94").Value = "07/10/2023 15:03"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/czech-republic/fnl/taborsko-sigma-olomouc/h80EI4z2/"

# Row 95
$ws.Range("F95").Value = "Vlasim"
$ws.Range("G95").Value = 3
$ws.Range("H95").Value = "Dukla Prague"
$ws.Range("I95").Value = 2
$ws.Range("J95").Value = 2.33
$ws.Range("K95").Value = "05/10/2023 08:12"
$ws.Range("L95").Value = 2.51
$ws.Range("M95").Value = "07/10/2023 15:57"
$ws.Range("N95").Value = 3.3
$ws.Range("O95").Value = "05/10/2023 08:12"
$ws.Range("P95").Value = 3.64
$ws.Range("Q95").Value = "07/10/2023 15:57"
$ws.Range("R95").Value = 2.7
$ws.Range("S95").Value = "05/10/2023 08:12"
$ws.Range("T95").Value = 2.6
$ws.Range("U95").Value = "07/10/2023 15:57"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/czech-republic/fnl/vlasim-dukla-prague/EH1AJpLe/"

# Rows 102, 103 (swap match data)
# Row 102
$ws.Range("F102").Value = "Lisen"
$ws.Range("G102").Value = 4
$ws.Range("H102").Value = "Sparta Prague B"
$ws.Range("I102").Value = 2
$ws.Range("J102").Value = 1.84
$ws.Range("K102").Value = "19/10/2023 09:12"
$ws.Range("L102").Value = 1.75
$ws.Range("M102").Value = "21/10/2023 15:24"
$ws.Range("N102").Value = 3.46
$ws.Range("O102").Value = "19/10/2023 09:12"
$ws.Range("P102").Value = 3.75
$ws.Range("Q102").Value = "21/10/2023 15:24"
$ws.Range("R102").Value = 3.73
$ws.Range("S102").Value = "19/10/2023 09:12"
$ws.Range("T102").Value = 4.59
$ws.Range("U102").Value = "21/10/2023 15:24"
$ws.Range("V102").Value = "https://www.betexplorer.com/football/czech-republic/fnl/lisen-sparta-prague/U7UeVsK7/"

# Row 103
$ws.Range("F103").Value = "Kromeriz"
$ws.Range("G103").Value = 1
$ws.Range("H103").Value = "Varnsdorf"
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 1.99
$ws.Range("K103").Value = "19/10/2023 09:12"
$ws.Range("L103").Value = 2.01
$ws.Range("M103").Value = "21/10/2023 15:21"
$ws.Range("N103").Value = 3.36
$ws.Range("O103").Value = "19/10/2023 09:12"
$ws.Range("P103").Value = 3.63
$ws.Range("Q103").Value = "21/10/2023 15:21"
$ws.Range("R103").Value = 3.32
$ws.Range("S103").Value = "19/10/2023 09:12"
$ws.Range("T103").Value = 3.52
$ws.Range("U103").Value = "21/10/2023 15:21"
$ws.Range("V103").Value = "https://www.betexplorer.com/football/czech-republic/fnl/kromeriz-varnsdorf/rBIwfRd8/"

# --- Append new match rows 106-113 (copy formatting from row 105 first) ---

# Row 106
$ws.Range("A105:V105").Copy()
$ws.Range("A106:V106").PasteSpecial(-4122)
$ws.Range("A106").Value = 105
$ws.Range("B106").Value = "czech-republic"
$ws.Range("C106").Value = "fnl"
$ws.Range("D106").Value = "2023-2024"
$ws.Range("E106").Value = 45226.70833333334
$ws.Range("F106").Value = "Brno"
$ws.Range("G106").Value = 3
$ws.Range("H106").Value = "Lisen"
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 1.68
$ws.Range("K106").Value = "26/10/2023 05:12"
$ws.Range("L106").Value = 1.78
$ws.Range("M106").Value = "27/10/2023 16:51"
$ws.Range("N106").Value = 3.63
$ws.Range("O106").Value = "26/10/2023 05:12"
$ws.Range("P106").Value = 3.47
$ws.Range("Q106").Value = "27/10/2023 16:51"
$ws.Range("R106").Value = 4.31
$ws.Range("S106").Value = "26/10/2023 05:12"
$ws.Range("T106").Value = 4.81
$ws.Range("U106").Value = "27/10/2023 16:51"
$ws.Range("V106").Value = "https://www.betexplorer.com/football/czech-republic/fnl/brno-lisen/8bkVqsZ0/"

# Row 107
$ws.Range("A105:V105").Copy()
$ws.Range("A107:V107").PasteSpecial(-4122)
$ws.Range("A107").Value = 106
$ws.Range("B107").Value = "czech-republic"
$ws.Range("C107").Value = "fnl"
$ws.Range("D107").Value = "2023-2024"
$ws.Range("E107").Value = 45226.75
$ws.Range("F107").Value = "Jihlava"
$ws.Range("G107").Value = 2
$ws.Range("H107").Value = "Pribram"
$ws.Range("I107").Value = 1
$ws.Range("J107").Value = 2.06
$ws.Range("K107").Value = "26/10/2023 06:12"
$ws.Range("L107").Value = 2.26
$ws.Range("M107").Value = "27/10/2023 17:03"
$ws.Range("N107").Value = 3.4
$ws.Range("O107").Value = "26/10/2023 06:12"
$ws.Range("P107").Value = 3.52
$ws.Range("Q107").Value = "27/10/2023 17:03"
$ws.Range("R107").Value = 3.1
$ws.Range("S107").Value = "26/10/2023 06:12"
$ws.Range("T107").Value = 3.02
$ws.Range("U107").Value = "27/10/2023 17:03"
$ws.Range("V107").Value = "https://www.betexplorer.com/football/czech-republic/fnl/jihlava-pribram/WWzQ5OeK/"

# Row 108
$ws.Range("A105:V105").Copy()
$ws.Range("A108:V108").PasteSpecial(-4122)
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "czech-republic"
$ws.Range("C108").Value = "fnl"
$ws.Range("D108").Value = "2023-2024"
$ws.Range("E108").Value = 45227.42708333334
$ws.Range("F108").Value = "Chrudim"
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = "Kromeriz"
$ws.Range("I108").Value = 1
$ws.Range("J108").Value = 1.93
$ws.Range("K108").Value = "26/10/2023 22:42"
$ws.Range("L108").Value = 1.87
$ws.Range("M108").Value = "28/10/2023 10:06"
$ws.Range("N108").Value = 3.41
$ws.Range("O108").Value = "26/10/2023 22:42"
$ws.Range("P108").Value = 3.48
$ws.Range("Q108").Value = "28/10/2023 10:12"
$ws.Range("R108").Value = 3.41
$ws.Range("S108").Value = "26/10/2023 22:42"
$ws.Range("T108").Value = 4.25
$ws.Range("U108").Value = "28/10/2023 10:06"
$ws.Range("V108").Value = "https://www.betexplorer.com/football/czech-republic/fnl/chrudim-kromeriz/U7jRpNJf/"

# Row 109
$ws.Range("A105:V105").Copy()
$ws.Range("A109:V109").PasteSpecial(-4122)
$ws.Range("A109").Value = 108
$ws.Range("B109").Value = "czech-republic"
$ws.Range("C109").Value = "fnl"
$ws.Range("D109").Value = "2023-2024"
$ws.Range("E109").Value = 45227.42708333334
$ws.Range("F109").Value = "Vyskov"
$ws.Range("G109").Value = 1
$ws.Range("H109").Value = "Opava"
$ws.Range("I109").Value = 3
$ws.Range("J109").Value = 1.9
$ws.Range("K109").Value = "26/10/2023 22:42"
$ws.Range("L109").Value = 1.87
$ws.Range("M109").Value = "28/10/2023 10:06"
$ws.Range("N109").Value = 3.71
$ws.Range("O109").Value = "26/10/2023 22:42"
$ws.Range("P109").Value = 3.46
$ws.Range("Q109").Value = "28/10/2023 10:06"
$ws.Range("R109").Value = 3.36
$ws.Range("S109").Value = "26/10/2023 22:42"
$ws.Range("T109").Value = 4.28
$ws.Range("U109").Value = "28/10/2023 10:06"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/czech-republic/fnl/mfk-vyskov-opava/GOhwrL3D/"

# Row 110
$ws.Range("A105:V105").Copy()
$ws.Range("A110:V110").PasteSpecial(-4122)
$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "czech-republic"
$ws.Range("C110").Value = "fnl"
$ws.Range("D110").Value = "2023-2024"
$ws.Range("E110").Value = 45227.4375
$ws.Range("F110").Value = "Prostejov"
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = "Dukla Prague"
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = 3.51
$ws.Range("K110").Value = "26/10/2023 22:42"
$ws.Range("L110").Value = 3.75
$ws.Range("M110").Value = "28/10/2023 10:17"
$ws.Range("N110").Value = 3.56
$ws.Range("O110").Value = "26/10/2023 22:42"
$ws.Range("P110").Value = 3.8
$ws.Range("Q110").Value = "28/10/2023 10:17"
$ws.Range("R110").Value = 1.86
$ws.Range("S110").Value = "26/10/2023 22:42"
$ws.Range("T110").Value = 1.89
$ws.Range("U110").Value = "28/10/2023 10:17"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/czech-republic/fnl/prostejov-dukla-prague/xYgZr1l7/"

# Row 111
$ws.Range("A105:V105").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)
$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "czech-republic"
$ws.Range("C111").Value = "fnl"
$ws.Range("D111").Value = "2023-2024"
$ws.Range("E111").Value = 45227.60416666666
$ws.Range("F111").Value = "Taborsko"
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = "Vlasim"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 2.02
$ws.Range("K111").Value = "27/10/2023 02:42"
$ws.Range("L111").Value = 1.79
$ws.Range("M111").Value = "28/10/2023 14:13"
$ws.Range("N111").Value = 3.5
$ws.Range("O111").Value = "27/10/2023 02:42"
$ws.Range("P111").Value = 3.87
$ws.Range("Q111").Value = "28/10/2023 14:22"
$ws.Range("R111").Value = 3.11
$ws.Range("S111").Value = "27/10/2023 02:42"
$ws.Range("T111").Value = 4.13
$ws.Range("U111").Value = "28/10/2023 14:21"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/czech-republic/fnl/taborsko-vlasim/z9cInqkr/"

# Row 112
$ws.Range("A105:V105").Copy()
$ws.Range("A112:V112").PasteSpecial(-4122)
$ws.Range("A112").Value = 111
$ws.Range("B112").Value = "czech-republic"
$ws.Range("C112").Value = "fnl"
$ws.Range("D112").Value = "2023-2024"
$ws.Range("E112").Value = 45228.4375
$ws.Range("F112").Value = "Sparta Prague B"
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = "Sigma Olomouc B"
$ws.Range("I112").Value = 2
$ws.Range("J112").Value = 2.03
$ws.Range("K112").Value = "27/10/2023 23:42"
$ws.Range("L112").Value = 2.19
$ws.Range("M112").Value = "29/10/2023 10:20"
$ws.Range("N112").Value = 3.47
$ws.Range("O112").Value = "27/10/2023 23:42"
$ws.Range("P112").Value = 3.42
$ws.Range("Q112").Value = "29/10/2023 10:20"
$ws.Range("R112").Value = 3.11
$ws.Range("S112").Value = "27/10/2023 23:42"
$ws.Range("T112").Value = 3.24
$ws.Range("U112").Value = "29/10/2023 10:20"
$ws.Range("V112").Value = "https://www.betexplorer.com/football/czech-republic/fnl/sparta-prague-sigma-olomouc/nq1Mo34l/"

# Row 113
$ws.Range("A105:V105").Copy()
$ws.Range("A113:V113").PasteSpecial(-4122)
$ws.Range("A113").Value = 112
$ws.Range("B113").Value = "czech-republic"
$ws.Range("C113").Value = "fnl"
$ws.Range("D113").Value = "2023-2024"
$ws.Range("E113").Value = 45228.60416666666
$ws.Range("F113").Value = "Varnsdorf"
$ws.Range("G113").Value = 4
$ws.Range("H113").Value = "Zizkov"
$ws.Range("I113").Value = 2
$ws.Range("J113").Value = 1.96
$ws.Range("K113").Value = "28/10/2023 03:43"
$ws.Range("L113").Value = 2.25
$ws.Range("M113").Value = "29/10/2023 14:20"
$ws.Range("N113").Value = 3.58
$ws.Range("O113").Value = "28/10/2023 03:43"
$ws.Range("P113").Value = 3.55
$ws.Range("Q113").Value = "29/10/2023 14:20"
$ws.Range("R113").Value = 3.19
$ws.Range("S113").Value = "28/10/2023 03:43"
$ws.Range("T113").Value = 3.01
$ws.Range("U113").Value = "29/10/2023 14:20"
$ws.Range("V113").Value = "https://www.betexplorer.com/football/czech-republic/fnl/varnsdorf-zizkov/ITWL64tE/"

$excel.Application.CutCopyMode = 0
